# "table of contents updated for PPT"
#
# 1. Refresh the cached "datetimeFigureOut" footer field (today's date,
#    dd-mm-yyyy) on the slide master and every slide layout:
#    24-08-2024 -> 02-09-2024
# 2. Rewrite the Contents (TOC) slide's bullet list with the new,
#    expanded table of contents.
# 3. Rename the "Identity Federation" section-title slide to
#    "Evolution of Identity Management" to match the new first TOC entry.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "24-08-2024") {
                $sh.TextFrame.TextRange.Text = "02-09-2024"
            }
        }
    }
}

# --- 1. Slide master + every custom (slide) layout -----------------------
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# --- 2. "Contents" slide (slide 2): rewrite the table of contents --------
$tocSlide = $p.Slides.Item(2)
$tocBody = $tocSlide.Shapes.Item(2)
$tocLines = @(
    "Evolution of Identity Management",
    "Core Concepts and Terminology",
    "Introduction to Identity Federation",
    "Identity Federation in the Cloud",
    "Identity Federation Protocols",
    "Federation Architecture",
    "Implementing Identity Federation",
    "Conclusion"
)
$tocBody.TextFrame.TextRange.Text = [string]::Join("`r", $tocLines)

# --- 3. "Identity Federation" slide (slide 3): retitle --------------------
$sectionSlide = $p.Slides.Item(3)
$sectionTitle = $sectionSlide.Shapes.Item(1)
$sectionTitle.TextFrame.TextRange.Text = "Evolution of Identity Management"
